# "Generate Report for Handoff"
#
# This regenerates the localization-status report data: a new file
# (4d338c82-aeeb-4631-9111-18f798c35898) has reached "Ready for handoff" and
# is inserted at the top of the list, and another new file
# (ae83df75-3af5-4da7-8e3c-3af509389c5f) has been appended at the end, in
# addition to the two pre-existing files. Overview/zh-cn/de-de sheets are all
# refreshed to reflect the new 4-row data set.

$wb = $excel.ActiveWorkbook

$mdBase  = "https://github.com/OpenLocalizationTest/oltest/blob/f04936a84c0a7bba850abe9fabdce3740b7bdc78/e2e"
$zhBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ea30fdd8fa5458b8607a290dd61582851074b99/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6cbbbdfb055c4099da7c6f0a7476e3669bd03185/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

# Ordered list of files now tracked by the report.
$files = @(
    @{ Id = "4d338c82-aeeb-4631-9111-18f798c35898"; Zh = "c8985388dad8fef4a9b0d57b007b35e6e8f5a357"; De = "c8985388dad8fef4a9b0d57b007b35e6e8f5a357" },
    @{ Id = "9df5e531-4c00-45b9-a522-98f78ba75baf"; Zh = "5d7cce0290ad35d371b0c8188d430edc00fbfdc3"; De = "5d7cce0290ad35d371b0c8188d430edc00fbfdc3" },
    @{ Id = "a9302f8a-bbe3-4db1-a066-cb2bd7605176"; Zh = "54e5b30ebe9f8c3b5bef3e33d29dd11a1d4180be"; De = "54e5b30ebe9f8c3b5bef3e33d29dd11a1d4180be" },
    @{ Id = "ae83df75-3af5-4da7-8e3c-3af509389c5f"; Zh = "3ed6f97c6450ef394e7658df06f9f6d352ec880a"; De = "3ed6f97c6450ef394e7658df06f9f6d352ec880a" }
)

$handoffDateOverview = "2016-21-13 22:21:09"
$handoffDatetimeZh = "2016-03-13 22:21:05"
$handoffDatetimeDe = "2016-03-13 22:21:09"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Hyperlinks.Delete()

for ($i = 0; $i -lt $files.Count; $i++) {
    $row = $i + 2
    $f = $files[$i]
    $mdName = "$($f.Id).md"

    $wsOverview.Range("B$row").Value = "Ready for handoff"
    $wsOverview.Range("C$row").Value = "Ready for handoff"
    $wsOverview.Range("D$row").Value = $handoffDateOverview

    $wsOverview.Hyperlinks.Add($wsOverview.Range("A$row"), "$mdBase/$mdName", [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
}

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item(2)
$wsZh.Hyperlinks.Delete()

for ($i = 0; $i -lt $files.Count; $i++) {
    $row = $i + 2
    $f = $files[$i]
    $mdName = "$($f.Id).md"
    $xlfName = "$($f.Id).$($f.Zh).zh-cn.xlf"

    $wsZh.Range("C$row").Value = "Ready for handoff"
    $wsZh.Range("E$row").Value = $handoffDatetimeZh
    $wsZh.Range("E$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $wsZh.Range("H$row").Value = "0001-01-01 00:00:00"
    $wsZh.Range("I$row").Value = "Include"

    $wsZh.Hyperlinks.Add($wsZh.Range("A$row"), "$mdBase/$mdName", [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
    $wsZh.Hyperlinks.Add($wsZh.Range("B$row"), "$mdBase/$mdName", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
    $wsZh.Hyperlinks.Add($wsZh.Range("D$row"), "$zhBase/$xlfName", [Type]::Missing, [Type]::Missing, $xlfName) | Out-Null
}

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item(3)
$wsDe.Hyperlinks.Delete()

for ($i = 0; $i -lt $files.Count; $i++) {
    $row = $i + 2
    $f = $files[$i]
    $mdName = "$($f.Id).md"
    $xlfName = "$($f.Id).$($f.De).de-de.xlf"

    $wsDe.Range("C$row").Value = "Ready for handoff"
    $wsDe.Range("E$row").Value = $handoffDatetimeDe
    $wsDe.Range("E$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $wsDe.Range("H$row").Value = "0001-01-01 00:00:00"
    $wsDe.Range("I$row").Value = "Include"

    $wsDe.Hyperlinks.Add($wsDe.Range("A$row"), "$mdBase/$mdName", [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
    $wsDe.Hyperlinks.Add($wsDe.Range("B$row"), "$mdBase/$mdName", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
    $wsDe.Hyperlinks.Add($wsDe.Range("D$row"), "$deBase/$xlfName", [Type]::Missing, [Type]::Missing, $xlfName) | Out-Null
}

Write-Output "Report regenerated for handoff"
